$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.058859333333333
$ws.Range("H2").Value = 6.176577999999999
$ws.Range("I2").Value = 0.03050820259458848
$ws.Range("J2").Value = 0.03050820259458848
$ws.Range("M2").Value = 11.146846
$ws.Range("N2").Value = 33.440538
$ws.Range("O2").Value = 0.2594806085672136
$ws.Range("P2").Value = 0.2594806085672136
$ws.Range("Q2").Value = 22.94978792432933
$ws.Range("R2").Value = 206.548091318964
$ws.Range("S2").Value = 0.007916286975535662
$ws.Range("T2").Value = 0.007916286975535662
$ws.Range("G3").Value = 2.058859333333333
$ws.Range("H3").Value = 6.176577999999999
$ws.Range("I3").Value = 0.03050820259458848
$ws.Range("J3").Value = 0.03050820259458848
$ws.Range("O3").Value = 0.6444737471070977
$ws.Range("P3").Value = 0.6444737471070977
$ws.Range("Q3").Value = 57.00054389642199
$ws.Range("R3").Value = 513.0048950677979
$ws.Range("S3").Value = 0.01966173564363692
$ws.Range("T3").Value = 0.01966173564363692
$ws.Range("G4").Value = 2.058859333333333
$ws.Range("H4").Value = 6.176577999999999
$ws.Range("I4").Value = 0.03050820259458848
$ws.Range("J4").Value = 0.03050820259458848
$ws.Range("O4").Value = 0.09604564432568881
$ws.Range("P4").Value = 0.09604564432568881
$ws.Range("Q4").Value = 8.494766450954888
$ws.Range("R4").Value = 76.452898058594
$ws.Range("S4").Value = 0.002930179975415902
$ws.Range("T4").Value = 0.002930179975415902
$ws.Range("I5").Value = 0.540047065760451
$ws.Range("J5").Value = 0.540047065760451
$ws.Range("M5").Value = 11.146846
$ws.Range("N5").Value = 33.440538
$ws.Range("O5").Value = 0.2594806085672136
$ws.Range("P5").Value = 0.2594806085672136
$ws.Range("Q5").Value = 406.2502728547213
$ws.Range("R5").Value = 3656.252455692492
$ws.Range("S5").Value = 0.1401317412784598
$ws.Range("T5").Value = 0.1401317412784598
$ws.Range("I6").Value = 0.540047065760451
$ws.Range("J6").Value = 0.540047065760451
$ws.Range("O6").Value = 0.6444737471070977
$ws.Range("P6").Value = 0.6444737471070977
$ws.Range("S6").Value = 0.3480461560848311
$ws.Range("T6").Value = 0.3480461560848311
$ws.Range("I7").Value = 0.540047065760451
$ws.Range("J7").Value = 0.540047065760451
$ws.Range("O7").Value = 0.09604564432568881
$ws.Range("P7").Value = 0.09604564432568881
$ws.Range("S7").Value = 0.05186916839716015
$ws.Range("T7").Value = 0.05186916839716015
$ws.Range("I8").Value = 0.4294447316449605
$ws.Range("J8").Value = 0.4294447316449605
$ws.Range("M8").Value = 11.146846
$ws.Range("N8").Value = 33.440538
$ws.Range("O8").Value = 0.2594806085672136
$ws.Range("P8").Value = 0.2594806085672136
$ws.Range("Q8").Value = 323.0496941245747
$ws.Range("R8").Value = 2907.447247121172
$ws.Range("S8").Value = 0.1114325803132181
$ws.Range("T8").Value = 0.1114325803132181
$ws.Range("I9").Value = 0.4294447316449605
$ws.Range("J9").Value = 0.4294447316449605
$ws.Range("O9").Value = 0.6444737471070977
$ws.Range("P9").Value = 0.6444737471070977
$ws.Range("S9").Value = 0.2767658553786297
$ws.Range("T9").Value = 0.2767658553786297
$ws.Range("I10").Value = 0.4294447316449605
$ws.Range("J10").Value = 0.4294447316449605
$ws.Range("O10").Value = 0.09604564432568881
$ws.Range("P10").Value = 0.09604564432568881
$ws.Range("S10").Value = 0.04124629595311275
$ws.Range("T10").Value = 0.04124629595311275
